$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for Camote (row 3 and row 4) had their values swapped:
# Row 3 (date, volume, min price, weighted avg price, price/kg) takes what was row 4's data
# Row 4 takes what was row 3's data.

# Row 3 new values (previously held by row 4)
$ws.Range("D3").Value = 44804
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 9500
$ws.Range("M3").Value = 9750
$ws.Range("P3").Value = 542

# Row 4 new values (previously held by row 3)
$ws.Range("D4").Value = 44714
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 9000
$ws.Range("M4").Value = 9500
$ws.Range("P4").Value = 528
